$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G2").Value = 500
$ws.Range("G3").Value = 3500
$ws.Range("G4").Value = 250
$ws.Range("G5").Value = 2000
$ws.Range("G6").Value = 12000
$ws.Range("G8").Value = 2000
$ws.Range("G9").Value = 1425
$ws.Range("G10").Value = 12000
$ws.Range("G11").Value = 3000
$ws.Range("G12").Value = 250
$ws.Range("G14").Value = 1500
$ws.Range("G15").Value = 1500
$ws.Range("G17").Value = 300
$ws.Range("G18").Value = 20000
$ws.Range("G19").Value = 60225
